$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 305378.25
$ws.Range("J17").Value = 347278.75
$ws.Range("L17").Value = 1041836.25
$ws.Range("N17").Value = -1042172.25
$ws.Range("H33").Value = 629.8125
$ws.Range("I33").Value = 321.625
$ws.Range("K33").Value = 321.625
$ws.Range("M33").Value = -92.625
$ws.Range("H53").Value = 1772.2
$ws.Range("I53").Value = 1687
$ws.Range("J53").Value = 1900
$ws.Range("K53").Value = 1687
$ws.Range("L53").Value = 1900
$ws.Range("M53").Value = -1050
$ws.Range("N53").Value = -3174
$ws.Range("H100").Value = 1025.3043
$ws.Range("I100").Value = 727.1177
$ws.Range("K100").Value = 727.1177
$ws.Range("M100").Value = -186.1177
$ws.Range("H111").Value = 3957.5
$ws.Range("I111").Value = 1861.3334
$ws.Range("J111").Value = 5529.625
$ws.Range("K111").Value = 5584.0002
$ws.Range("L111").Value = 16588.875
$ws.Range("M111").Value = -2517.0002
$ws.Range("N111").Value = -22722.875
$ws.Range("H129").Value = 2391.2856
$ws.Range("I129").Value = 1996
$ws.Range("K129").Value = 5988
$ws.Range("M129").Value = -988
$ws.Range("H132").Value = 1367.695
$ws.Range("I132").Value = 1053.6364
$ws.Range("J132").Value = 2288.9333
$ws.Range("K132").Value = 3160.9092
$ws.Range("L132").Value = 6866.7999
$ws.Range("M132").Value = -630.9092000000001
$ws.Range("N132").Value = -11926.7999
$ws.Range("H137").Value = 3891028.5
$ws.Range("I137").Value = 5057278
$ws.Range("K137").Value = 15171834
$ws.Range("M137").Value = -15169284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4629.1
$ws.Range("I61").Value = 2134.9092
$ws.Range("K61").Value = 2134.9092
$ws.Range("M61").Value = -1922.9092
$ws.Range("H102").Value = 1869.2
$ws.Range("I102").Value = 1684.4445
$ws.Range("K102").Value = 1684.4445
$ws.Range("M102").Value = -62.44450000000006
$ws.Range("H110").Value = 8481.925999999999
$ws.Range("I110").Value = 7196.2173
$ws.Range("K110").Value = 7196.2173
$ws.Range("M110").Value = -5151.2173
$ws.Range("H132").Value = 2915.1875
$ws.Range("I132").Value = 3597.7827
$ws.Range("K132").Value = 10793.3481
$ws.Range("M132").Value = -8263.348100000001
$ws.Range("H136").Value = 4629.1
$ws.Range("I136").Value = 2134.9092
$ws.Range("K136").Value = 6404.7276
$ws.Range("M136").Value = -3854.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3740.9714
$ws.Range("I99").Value = 3255
$ws.Range("J99").Value = 4469.9287
$ws.Range("K99").Value = 3255
$ws.Range("L99").Value = 4469.9287
$ws.Range("M99").Value = -1757
$ws.Range("N99").Value = -7465.9287
$ws.Range("H105").Value = 2158.3635
$ws.Range("I105").Value = 2158.3635
$ws.Range("K105").Value = 2158.3635
$ws.Range("M105").Value = -411.3634999999999
$ws.Range("H107").Value = 913.96155
$ws.Range("I107").Value = 1098.5
$ws.Range("K107").Value = 1098.5
$ws.Range("M107").Value = 821.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2243.5715
$ws.Range("I16").Value = 2243.5715
$ws.Range("K16").Value = 2243.5715
$ws.Range("M16").Value = -1956.5715
$ws.Range("H31").Value = 201626.67
$ws.Range("I31").Value = 279828.9
$ws.Range("J31").Value = 45222.168
$ws.Range("K31").Value = 279828.9
$ws.Range("L31").Value = 45222.168
$ws.Range("M31").Value = -279533.9
$ws.Range("N31").Value = -45812.168
$ws.Range("H34").Value = 201626.67
$ws.Range("I34").Value = 279828.9
$ws.Range("J34").Value = 45222.168
$ws.Range("K34").Value = 279828.9
$ws.Range("L34").Value = 45222.168
$ws.Range("M34").Value = -279626.9
$ws.Range("N34").Value = -45626.168
$ws.Range("H113").Value = 2243.5715
$ws.Range("I113").Value = 2243.5715
$ws.Range("K113").Value = 2243.5715
$ws.Range("M113").Value = -73.57150000000001
$ws.Range("H132").Value = 2699.5625
$ws.Range("I132").Value = 2746.8667
$ws.Range("K132").Value = 8240.6001
$ws.Range("M132").Value = -5710.6001
$ws.Range("H134").Value = 5073.9355
$ws.Range("I134").Value = 5783.25
$ws.Range("K134").Value = 17349.75
$ws.Range("M134").Value = -14814.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1749.6666
$ws.Range("I11").Value = 1749.6666
$ws.Range("K11").Value = 5248.9998
$ws.Range("M11").Value = -5108.9998
$ws.Range("H80").Value = 1899
$ws.Range("I80").Value = 1898.5
$ws.Range("K80").Value = 5695.5
$ws.Range("M80").Value = -4759.5
$ws.Range("H83").Value = 1899
$ws.Range("I83").Value = 1898.5
$ws.Range("K83").Value = 17086.5
$ws.Range("M83").Value = -12406.5
$ws.Range("H92").Value = 495.2857
$ws.Range("I92").Value = 516.75
$ws.Range("J92").Value = 466.66666
$ws.Range("K92").Value = 1550.25
$ws.Range("L92").Value = 1399.99998
$ws.Range("M92").Value = -302.25
$ws.Range("N92").Value = -3895.99998
$ws.Range("H129").Value = 1880.3889
$ws.Range("I129").Value = 1560.875
$ws.Range("J129").Value = 2136
$ws.Range("K129").Value = 4682.625
$ws.Range("L129").Value = 6408
$ws.Range("M129").Value = 317.375
$ws.Range("N129").Value = -16408
$ws.Range("H131").Value = 1435.4458
$ws.Range("I131").Value = 1009.6
$ws.Range("K131").Value = 3028.8
$ws.Range("M131").Value = 2011.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7414.875
$ws.Range("I70").Value = 7916.0586
$ws.Range("J70").Value = 6197.7144
$ws.Range("K70").Value = 7916.0586
$ws.Range("L70").Value = 6197.7144
$ws.Range("M70").Value = -7646.0586
$ws.Range("N70").Value = -6737.7144
$ws.Range("H73").Value = 7414.875
$ws.Range("I73").Value = 7916.0586
$ws.Range("J73").Value = 6197.7144
$ws.Range("K73").Value = 7916.0586
$ws.Range("L73").Value = 6197.7144
$ws.Range("M73").Value = -6980.0586
$ws.Range("N73").Value = -8069.7144
$ws.Range("H113").Value = 3189.7407
$ws.Range("J113").Value = 4331.375
$ws.Range("L113").Value = 4331.375
$ws.Range("N113").Value = -8671.375
$ws.Range("H123").Value = 21700.125
$ws.Range("J123").Value = 21700.125
$ws.Range("L123").Value = 21700.125
$ws.Range("N123").Value = -26600.125
$ws.Range("H132").Value = 46575.656
$ws.Range("I132").Value = 56848
$ws.Range("K132").Value = 170544
$ws.Range("M132").Value = -168014

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8145.25
$ws.Range("I7").Value = 8203.138999999999
$ws.Range("K7").Value = 8203.138999999999
$ws.Range("M7").Value = -8091.138999999999
$ws.Range("H22").Value = 1382.85
$ws.Range("J22").Value = 1393.375
$ws.Range("L22").Value = 1393.375
$ws.Range("N22").Value = -1983.375
$ws.Range("H27").Value = 1382.85
$ws.Range("J27").Value = 1393.375
$ws.Range("L27").Value = 1393.375
$ws.Range("N27").Value = -1607.375
$ws.Range("H68").Value = 5505.3887
$ws.Range("I68").Value = 6035.643
$ws.Range("J68").Value = 3649.5
$ws.Range("K68").Value = 6035.643
$ws.Range("L68").Value = 3649.5
$ws.Range("M68").Value = -5286.643
$ws.Range("N68").Value = -5147.5
$ws.Range("H71").Value = 5505.3887
$ws.Range("I71").Value = 6035.643
$ws.Range("J71").Value = 3649.5
$ws.Range("K71").Value = 30178.215
$ws.Range("L71").Value = 18247.5
$ws.Range("M71").Value = -26434.215
$ws.Range("N71").Value = -25735.5
$ws.Range("H93").Value = 2743.8667
$ws.Range("I93").Value = 988
$ws.Range("J93").Value = 5377.6665
$ws.Range("K93").Value = 988
$ws.Range("L93").Value = 5377.6665
$ws.Range("M93").Value = 260
$ws.Range("N93").Value = -7873.6665
$ws.Range("H122").Value = 7852
$ws.Range("I122").Value = 7478.2856
$ws.Range("J122").Value = 8898.4
$ws.Range("K122").Value = 22434.8568
$ws.Range("L122").Value = 26695.2
$ws.Range("M122").Value = -19984.8568
$ws.Range("N122").Value = -31595.2
$ws.Range("H126").Value = 8145.25
$ws.Range("I126").Value = 8203.138999999999
$ws.Range("K126").Value = 24609.417
$ws.Range("M126").Value = -22139.417
$ws.Range("H132").Value = 9946.75
$ws.Range("I132").Value = 9898.75
$ws.Range("J132").Value = 9994.75
$ws.Range("K132").Value = 29696.25
$ws.Range("L132").Value = 29984.25
$ws.Range("M132").Value = -27166.25
$ws.Range("N132").Value = -35044.25
$ws.Range("H134").Value = 66406.5
$ws.Range("J134").Value = 78423
$ws.Range("L134").Value = 78423
$ws.Range("N134").Value = -88563
$ws.Range("H141").Value = 59816.668
$ws.Range("J141").Value = 59816.668
$ws.Range("L141").Value = 59816.668
$ws.Range("N141").Value = -70176.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 857.2857
$ws.Range("I100").Value = 723.875
$ws.Range("K100").Value = 1447.75
$ws.Range("M100").Value = -906.75
$ws.Range("H122").Value = 2901.4
$ws.Range("I122").Value = 2901.4
$ws.Range("K122").Value = 8704.200000000001
$ws.Range("M122").Value = -6254.200000000001
$ws.Range("H132").Value = 1164.7273
$ws.Range("I132").Value = 1201.5555
$ws.Range("K132").Value = 3604.6665
$ws.Range("M132").Value = -1074.6665

Write-Output "edit complete"